$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Gender_ID"
$ws.Range("C2").Value = "Name"
